$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer")
$ws.Activate()

# Update the formulas in B2:B7 - remove the CONCATENATE("return """..."""";") wrapper
# keeping just the inner SUBSTITUTE(UPPER(LEFT(...)),...) expression.
$ws.Range("B2").Formula = '=SUBSTITUTE(UPPER(LEFT(A2,4)), " ", "")'
$ws.Range("B3").Formula = '=SUBSTITUTE(UPPER(LEFT(A3,4)), " ", "")'
$ws.Range("B4").Formula = '=SUBSTITUTE(UPPER(LEFT(A4,4)), " ", "")'
$ws.Range("B5").Formula = '=SUBSTITUTE(UPPER(LEFT(A5,4)), " ", "")'
$ws.Range("B6").Formula = '=SUBSTITUTE(UPPER(LEFT(A6,4)), " ", "")'
$ws.Range("B7").Formula = '=SUBSTITUTE(UPPER(LEFT(A7,4)), " ", "")'

$excel.Calculate()

# Update the selection on the active sheet to B2:B7, active cell B2.
$ws.Range("B2:B7").Select()
